$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$d = (Get-Date -Year 2015 -Month 5 -Day 13).Date
$ws.Range("P5").Value = $d
$ws.Range("P6").Value = $d
$ws.Range("P7").Value = $d
$ws.Range("P8").Value = $d

$ws.Range("AN6").Value = "OPEN_ACCESS"
$ws.Range("AN7").Value = "OPEN_ACCESS"
$ws.Range("AN8").Value = "OPEN_ACCESS"
$ws.Range("AN9").Value = "OPEN_ACCESS"

$ws.Range("AP5").Select() | Out-Null
